$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a date-shaped literal string into a cell without Excel's
# automatic "looks like a date" input parsing turning it into a real date
# (and thereby introducing a date number-format style). We build the text
# via a formula (so it is produced as a string, not auto-coerced), copy it,
# and paste-special "values only" into the destination - pasted text values
# are not re-parsed as dates, so the destination keeps the default style.
$scratch = $ws.Range("Z100")

function Set-LiteralText {
    param($targetCell, [string]$text)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)
}

# Row 3 updates
$ws.Range("D3").Value = 11111.0
$ws.Range("C3").Value = "donat"
Set-LiteralText $ws.Range("F3") "2024-11-23"

# Row 4 (new row)
$ws.Range("A4").Value = 102.0
$ws.Range("B4").Value = "anette"
$ws.Range("C4").Value = "donat"
$ws.Range("D4").Value = 22333.0
Set-LiteralText $ws.Range("E4") "2024-11-18"
Set-LiteralText $ws.Range("F4") "2024-11-23"

# Clean up the scratch cell/clipboard so nothing extra leaks into the sheet.
$scratch.Clear()
$excel.CutCopyMode = $false
